# Auto-generated edit script applying crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: D2: '87.193.02'->'87.193.43'; E2: '  -3.07%  '->'  -3.23%  '
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '87.193.43'
$ws.Range("E2").Value = '  -3.23%  '

# Row 3: D3: '3.022.41'->'3.020.05'; E3: '  -6.49%  '->'  -6.66%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.020.05'
$ws.Range("E3").Value = '  -6.66%  '

# Row 4: E4: '  +0.11%  '->'  -0.02%  '
$ws.Range("E4").Value = '  -0.02%  '

# Row 5: D5: '205.04'->'205.31'; E5: '  -6.37%  '->'  -6.67%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '205.31'
$ws.Range("E5").Value = '  -6.67%  '

# Row 6: D6: '608.68'->'609.54'; E6: '  -3.44%  '->'  -3.62%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '609.54'
$ws.Range("E6").Value = '  -3.62%  '

# Row 7: E7: '  -9.47%  '->'  -9.31%  '
$ws.Range("E7").Value = '  -9.31%  '

# Row 8: D8: '0.800'->'0.802'; E8: '  +14.19%  '->'  +14.37%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.802'
$ws.Range("E8").Value = '  +14.37%  '

# Row 9: E9: '  +0.13%  '->'  +0.04%  '
$ws.Range("E9").Value = '  +0.04%  '

# Row 10: D10: '3.021.98'->'3.018.60'; E10: '  -6.44%  '->'  -6.60%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.018.60'
$ws.Range("E10").Value = '  -6.60%  '

# Row 11: D11: '0.579'->'0.583'; E11: '  +1.05%  '->'  +1.65%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.583'
$ws.Range("E11").Value = '  +1.65%  '

# Row 12: D12: '0.177'->'0.176'; E12: '  -1.08%  '->'  -1.32%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.176'
$ws.Range("E12").Value = '  -1.32%  '

# Row 13: D13: '0.0000228'->'0.0000227'; E13: '  -12.13%  '->'  -12.76%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000227'
$ws.Range("E13").Value = '  -12.76%  '

# Row 14: D14: '5.16'->'5.17'; E14: '  -4.35%  '->'  -4.25%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.17'
$ws.Range("E14").Value = '  -4.25%  '

# Row 15: D15: '87.239.91'->'87.015.45'; E15: '  -2.54%  '->'  -2.93%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '87.015.45'
$ws.Range("E15").Value = '  -2.93%  '

# Row 16: D16: '3.589.36'->'3.579.24'; E16: '  -5.81%  '->'  -6.47%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.579.24'
$ws.Range("E16").Value = '  -6.47%  '

# Row 17: D17: '30.74'->'30.77'; E17: '  -7.92%  '->'  -8.09%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.77'
$ws.Range("E17").Value = '  -8.09%  '

# Row 18: D18: '3.047.84'->'3.038.20'; E18: '  -4.82%  '->'  -5.32%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.038.20'
$ws.Range("E18").Value = '  -5.32%  '

# Row 19: D19: '3.11'->'3.12'; E19: '  -7.06%  '->'  -6.86%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.12'
$ws.Range("E19").Value = '  -6.86%  '

# Row 20: E20: '  -18.22%  '->'  -18.38%  '
$ws.Range("E20").Value = '  -18.38%  '

# Row 21: D21: '12.77'->'12.79'; E21: '  -5.55%  '->'  -6.23%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.79'
$ws.Range("E21").Value = '  -6.23%  '

# Row 22: D22: '411.57'->'411.46'; E22: '  -5.94%  '->'  -6.33%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '411.46'
$ws.Range("E22").Value = '  -6.33%  '

# Row 23: D23: '7.90'->'7.89'; E23: '  -8.60%  '->'  -8.82%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.89'
$ws.Range("E23").Value = '  -8.82%  '

# Row 24: E24: '  -7.31%  '->'  -7.67%  '
$ws.Range("E24").Value = '  -7.67%  '

# Row 25: D25: '5.22'->'5.21'; E25: '  +0.92%  '->'  +0.29%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.21'
$ws.Range("E25").Value = '  +0.29%  '

# Row 26: B26: 'Litecoin'->'Aptos'; C26: 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'->'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; D26: '79.61'->'11.25'; E26: '  +0.31%  '->'  -5.38%  '
$ws.Range("B26").Value = 'Aptos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.25'
$ws.Range("E26").Value = '  -5.38%  '

# Row 27: B27: 'Aptos'->'Litecoin'; C27: 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'->'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; D27: '11.19'->'79.42'; E27: '  -5.50%  '->'  -0.26%  '
$ws.Range("B27").Value = 'Litecoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '79.42'
$ws.Range("E27").Value = '  -0.26%  '

# Row 28: D28: '3.227.18'->'3.225.30'; E28: '  -4.56%  '->'  -4.64%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.225.30'
$ws.Range("E28").Value = '  -4.64%  '

# Row 29: D29: '0.999'->'1.00'; E29: '  -0.07%  '->'  -0.01%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  -0.01%  '

# Row 30: D30: '1.07'->'1.08'; E30: '  +7.27%  '->'  +7.90%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.08'
$ws.Range("E30").Value = '  +7.90%  '

# Row 31: D31: '0.154'->'0.155'; E31: '  -2.74%  '->'  -1.85%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.155'
$ws.Range("E31").Value = '  -1.85%  '

# Row 32: D32: '7.87'->'7.85'; E32: '  -7.73%  '->'  -8.42%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.85'
$ws.Range("E32").Value = '  -8.42%  '

# Row 33: D33: '492.26'->'492.20'; E33: '  -9.33%  '->'  -9.54%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '492.20'
$ws.Range("E33").Value = '  -9.54%  '

# Row 34: D34: '3.39'->'3.40'; E34: '  -17.98%  '->'  -17.29%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.40'
$ws.Range("E34").Value = '  -17.29%  '

# Row 35: E35: '  -8.19%  '->'  -8.28%  '
$ws.Range("E35").Value = '  -8.28%  '

# Row 36: D36: '6.40'->'6.39'; E36: '  -8.84%  '->'  -8.99%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.39'
$ws.Range("E36").Value = '  -8.99%  '

# Row 37: D37: '1.20'->'1.19'; E37: '  -8.07%  '->'  -8.32%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.19'
$ws.Range("E37").Value = '  -8.32%  '

# Row 38: E38: '  -0.99%  '->'  -1.08%  '
$ws.Range("E38").Value = '  -1.08%  '

# Row 39: D39: '21.69'->'21.68'; E39: '  -3.57%  '->'  -3.73%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '21.68'
$ws.Range("E39").Value = '  -3.73%  '

# Row 40: E40: '  -0.45%  '->'  -0.68%  '
$ws.Range("E40").Value = '  -0.68%  '

# Row 41: E41: '  +0.37%  '->'  +0.31%  '
$ws.Range("E41").Value = '  +0.31%  '

# Row 42: E42: '  -0.04%  '->'  -0.07%  '
$ws.Range("E42").Value = '  -0.07%  '

# Row 43: D43: '148.12'->'148.15'; E43: '  +0.22%  '->'  +0.27%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '148.15'
$ws.Range("E43").Value = '  +0.27%  '

# Row 44: E44: '  -6.57%  '->'  -6.77%  '
$ws.Range("E44").Value = '  -6.77%  '

# Row 45: E45: '  +6.08%  '->'  +6.04%  '
$ws.Range("E45").Value = '  +6.04%  '

# Row 46: E46: '  -9.06%  '->'  -9.49%  '
$ws.Range("E46").Value = '  -9.49%  '

# Row 47: D47: '43.15'->'43.09'; E47: '  -1.30%  '->'  -1.63%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '43.09'
$ws.Range("E47").Value = '  -1.63%  '

# Row 48: D48: '0.0650'->'0.0649'; E48: '  +7.05%  '->'  +6.45%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0649'
$ws.Range("E48").Value = '  +6.45%  '

# Row 49: B49: 'Aave'->'Mantle'; C49: 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'->'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'; D49: '152.05'->'0.682'; E49: '  -12.39%  '->'  -9.66%  '
$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.682'
$ws.Range("E49").Value = '  -9.66%  '

# Row 50: B50: 'Mantle'->'Aave'; C50: 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'->'https://coinranking.com/coin/ixgUfzmLR+aave-aave'; D50: '0.683'->'151.31'; E50: '  -9.04%  '->'  -13.13%  '
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '151.31'
$ws.Range("E50").Value = '  -13.13%  '

# Row 51: E51: '  -9.31%  '->'  -9.38%  '
$ws.Range("E51").Value = '  -9.38%  '
